$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an unmodified cell style reference (row 1 header cells are untouched,
# but we need a plain/default-style cell from column D/E; row bodies have no explicit style)
$defaultStyleCell = $ws.Range("B2")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.417.73"
$ws.Range("D2").Style = $defaultStyleCell.Style
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.05"
$ws.Range("D3").Style = $defaultStyleCell.Style
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  -0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.21"
$ws.Range("D5").Style = $defaultStyleCell.Style
$ws.Range("E5").Value = "  -2.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("D6").Style = $defaultStyleCell.Style
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4451"
$ws.Range("D7").Style = $defaultStyleCell.Style
$ws.Range("E7").Value = "  +5.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3771"
$ws.Range("D8").Style = $defaultStyleCell.Style
$ws.Range("E8").Value = "  +7.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.82"
$ws.Range("D9").Style = $defaultStyleCell.Style
$ws.Range("E9").Value = "  -1.63%  "

$ws.Range("E10").Value = "  +0.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07501"
$ws.Range("D11").Style = $defaultStyleCell.Style
$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.57"
$ws.Range("D12").Style = $defaultStyleCell.Style
$ws.Range("E12").Value = "  -1.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.000"
$ws.Range("D13").Style = $defaultStyleCell.Style
$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.617"
$ws.Range("D14").Style = $defaultStyleCell.Style
$ws.Range("E14").Value = "  +4.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.297"
$ws.Range("D15").Style = $defaultStyleCell.Style
$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.802.31"
$ws.Range("D16").Style = $defaultStyleCell.Style
$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001092"
$ws.Range("D17").Style = $defaultStyleCell.Style
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06807"
$ws.Range("D18").Style = $defaultStyleCell.Style
$ws.Range("E18").Value = "  +1.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.57"
$ws.Range("D19").Style = $defaultStyleCell.Style
$ws.Range("E19").Value = "  -2.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9992"
$ws.Range("D20").Style = $defaultStyleCell.Style
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.59"
$ws.Range("D21").Style = $defaultStyleCell.Style
$ws.Range("E21").Value = "  +1.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.326"
$ws.Range("D22").Style = $defaultStyleCell.Style
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.414.13"
$ws.Range("D23").Style = $defaultStyleCell.Style
$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.81"
$ws.Range("D24").Style = $defaultStyleCell.Style
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.413"
$ws.Range("D25").Style = $defaultStyleCell.Style
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.47"
$ws.Range("D26").Style = $defaultStyleCell.Style
$ws.Range("E26").Value = "  -1.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.99"
$ws.Range("D27").Style = $defaultStyleCell.Style
$ws.Range("E27").Value = "  -1.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.350"
$ws.Range("D28").Style = $defaultStyleCell.Style
$ws.Range("E28").Value = "  -4.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.006.11"
$ws.Range("D29").Style = $defaultStyleCell.Style
$ws.Range("E29").Value = "  -0.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.34"
$ws.Range("D30").Style = $defaultStyleCell.Style
$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("E31").Value = "  -3.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.005"
$ws.Range("D32").Style = $defaultStyleCell.Style
$ws.Range("E32").Value = "  -1.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.817"
$ws.Range("D33").Style = $defaultStyleCell.Style
$ws.Range("E33").Value = "  -2.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09339"
$ws.Range("D34").Style = $defaultStyleCell.Style
$ws.Range("E34").Value = "  +1.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2277"
$ws.Range("D35").Style = $defaultStyleCell.Style
$ws.Range("E35").Value = "  +4.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.13"
$ws.Range("D36").Style = $defaultStyleCell.Style
$ws.Range("E36").Value = "  -1.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06352"
$ws.Range("D37").Style = $defaultStyleCell.Style
$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02343"
$ws.Range("D38").Style = $defaultStyleCell.Style
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6592"
$ws.Range("D39").Style = $defaultStyleCell.Style
$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.158"
$ws.Range("D40").Style = $defaultStyleCell.Style
$ws.Range("E40").Value = "  -1.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.209"
$ws.Range("D41").Style = $defaultStyleCell.Style
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.453"
$ws.Range("D42").Style = $defaultStyleCell.Style
$ws.Range("E42").Value = "  -3.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.089"
$ws.Range("D43").Style = $defaultStyleCell.Style
$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9986"
$ws.Range("D44").Style = $defaultStyleCell.Style
$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.81"
$ws.Range("D45").Style = $defaultStyleCell.Style
$ws.Range("E45").Value = "  -3.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6083"
$ws.Range("D46").Style = $defaultStyleCell.Style
$ws.Range("E46").Value = "  -0.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.814"
$ws.Range("D47").Style = $defaultStyleCell.Style
$ws.Range("E47").Value = "  -1.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.42"
$ws.Range("D48").Style = $defaultStyleCell.Style
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.030"
$ws.Range("D49").Style = $defaultStyleCell.Style
$ws.Range("E49").Value = "  -1.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07088"
$ws.Range("D50").Style = $defaultStyleCell.Style
$ws.Range("E50").Value = "  -0.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.155"
$ws.Range("D51").Style = $defaultStyleCell.Style
$ws.Range("E51").Value = "  -1.97%  "
